$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.807.27"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.631.37"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.14"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.64"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "1.855.18"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "1.605.59"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.70"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "25.798.09"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.998"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.95"
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.30"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.80"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.22"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.904"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "1.144.39"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.544"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.38"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.802"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "1.765.27"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.44"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.46"
$ws.Range("E49").Value = "  +5.51%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.57"
$ws.Range("E51").Value = "  +0.11%  "
